# This workbook lists weekly Jengibre (ginger) price observations.
# The edit re-shuffles the per-row observations (Fecha/Volumen/Precio
# minimo/maximo/promedio ponderado/Precio $/Kg) among the existing rows,
# while rows 5, 10 and 15 stay untouched.
#
# Strategy: snapshot the "before" values for columns D, J, K, L, M, P for
# every affected row, then re-assign them to their new rows according to
# the row -> source-row mapping derived from the diff. Reading everything
# up front avoids clobbering source data before it has been copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the shuffle.
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Rows that participate in the shuffle (5, 10, 15 are left as-is).
$rows = @(2, 3, 4, 6, 7, 8, 9, 11, 12, 13, 14, 16, 17)

# Snapshot current ("before") values for each involved row/column.
# NOTE: use Value2 (not Value) when reading through this COM shim --
# Value is exposed as a parameterized property and round-tripping it
# through a PowerShell variable does not yield the underlying scalar.
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Target row <- source row mapping (which row's old data becomes the
# new data for the given row).
$mapping = @{
    2  = 13
    3  = 6
    4  = 9
    6  = 8
    7  = 17
    8  = 12
    9  = 4
    11 = 7
    12 = 14
    13 = 2
    14 = 16
    16 = 11
    17 = 3
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $sourceVals[$c]
    }
}
